$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Developed"/"Emerging" quartile blocks each shrink from 5 columns
# (quartiles 0-4, i.e. B:F and G:K) to 4 columns (quartiles 0-3, i.e. B:E
# and F:I). Deleting one column out of each block reproduces this layout
# while letting Excel's own column-delete logic take care of shifting the
# remaining cells/merged ranges/labels left, instead of manually tearing
# down and rebuilding the merges (which would also re-flow border
# formatting unnecessarily).
#
# Deleting column F removes the trailing "4" quartile under "Developed",
# shrinks the B1:F1 merge down to B1:E1, and shifts the G1:K1 ("Emerging")
# merge left to F1:J1 - sliding the "Emerging" label into its new anchor
# cell F1 for free.
$ws.Range("F1").EntireColumn.Delete()

# Deleting the (now) last column J removes the trailing "4" quartile under
# "Emerging" and shrinks the F1:J1 merge down to F1:I1.
$ws.Range("J1").EntireColumn.Delete()

# --- Row 4: refresh with the newly processed modeling data values.
$ws.Range("B4").Value = 0.008051437324986817
$ws.Range("C4").Value = 0.006563553620001397
$ws.Range("D4").Value = 0.007846016949009856
$ws.Range("E4").Value = 0.006819335881800336
$ws.Range("F4").Value = 0.01386722714049879
$ws.Range("G4").Value = 0.01316462850688125
$ws.Range("H4").Value = 0.008913264273177205
$ws.Range("I4").Value = 0.01287100601400059
